# Weekly data refresh: insert the newest week's price record for Mango at
# "Vega Modelo de Temuco" right before the existing row 159, pushing the
# rest of that sub-block (old rows 159-213) down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 159:213 down to 160:214 (inherits formatting from the row
# above, matching the existing date-format style on column D).
$ws.Rows.Item(159).Insert()

# Populate the newly inserted row 159 with this week's record.
$ws.Range("A159").Value = 10
$ws.Range("B159").Value = "Vega Modelo de Temuco"
$ws.Range("C159").Value = "La Araucanía"
$ws.Range("D159").Value = 44468
$ws.Range("E159").Value = 9
$ws.Range("F159").Value = "Fruta"
$ws.Range("G159").Value = 100108
$ws.Range("H159").Value = "Tropicales y subtropicales"
$ws.Range("I159").Value = 100108002
$ws.Range("J159").Value = "Mango"
$ws.Range("K159").Value = "Sin especificar"
$ws.Range("L159").Value = "Primera"
$ws.Range("M159").Value = 300
$ws.Range("N159").Value = 9000
$ws.Range("O159").Value = 9000
$ws.Range("P159").Value = 9000
$ws.Range("Q159").Value = "$/bandeja 4 kilos"
$ws.Range("R159").Value = "Brasil"
$ws.Range("S159").Value = 2250
$ws.Range("T159").Value = 4
